$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 71: a "skipped day" row (Ser + Date only, like row 68) ---
$ws.Rows.Item(70).Copy()
$ws.Rows.Item(71).Insert(-4121)

# --- Row 72: the new article row (full content, like rows 69/70) ---
$ws.Rows.Item(70).Copy()
$ws.Rows.Item(72).Insert(-4121)

# Row 71 values
$ws.Range("A71").Value = 70
$ws.Range("B71").Value = 43900
$ws.Range("C71:F71").Clear()

# Row 72 values
$ws.Range("A72").Value = 71
$ws.Range("B72").Value = 43901
$ws.Range("C72").Value = "Surah Al Maida, verse 101 - 103"
$ws.Range("F72").Value = "Importance of fighting your fears, Results of a worried life, Building self confidence, Try to relax in a busy day, Balloons"
$ws.Range("D72").Value = @'
h1: 3 Things to make your life easier 
p.note: I am only 32 years. I am not sure if this is the ideal formula to live a balanced life. But I have seen many things getting fixed real quick when I banked on these 3 habits. 😇
img: struggle.jpg
p.caption.text-center: Photo by <a href=https://unsplash.com/@dj_johns1 target=_blank> DJ Johnson</a> on Unsplash
p: This is what all great people tell. Trust your steps. I am may be not very experienced and not very mature. But I am moving at a lightening speed to hit the sweet spot of victory. It is feeling great. Not because of all the successes but because of all the blessings in every day failures.<b> Some day it is a blog not getting good audience. Other days it is a meetup not getting good feedback. </b>In spite of major failures and set backs every month, life is calm, sweet and bright. 
p: I will not ask you to write blogs. Neither I will ask you to learn to code and create websites. Because all these skills and habits are a by product of a focused life. 
p: Three things 👇
h3: <span class=lavendar>1. </span> Stop worrying about Tomorrow, Money and Kids
img: future.jpg
p.caption.text-center: Photo by <a href=https://unsplash.com/@juvx target=_blank> Vlad Busuioc </a> on Unsplash
p: Give up on all these parts of your life. Do not try to control your future. Do not try to accumulate wealth. And do not decide your kids’ future. Give up on your plans. Stop thinking about tomorrow. Close this book. Place it in that highest shelf on your rack. Make it unreachable. 
p: More thinking about future is equivalent of creating more doors to your life. You will need to keep an eye on each of these doors to make sure dullness, sadness and failures do not hit you down. Create lesser doors. One if possible. That is focused only on your own self. That tells you to improve you only. Asks you to be childish, selfless and not-worried. Block all other doors. Shut them down. Brick them up. 
h3: <span class=lavendar>2. </span> Trust yourself to create the new Normal
img: confidence.jpg
p.caption.text-center: Photo by <a href=https://unsplash.com/@webdonut target=_blank> Web Donut </a>on Unsplash
p: Do not under estimate your own potential. Do not consider yourself any lesser than Von Gough or Iqbal. You are more valueable, pricely and beautiful than all other people around you. Each step you take now, has more potential of creating dents than all great leaders in the history. 
p: You have technology in your fingers’ reach. One word out of your hand has more reach than all the words combined until yesterday. One truth spoken today, will create more impact than all the truths till yesterday. You are valuable, beautiful and talented. 
p: Create your new Normal. Take the initiative. Move boldly. Loudly. Confidently. Life is slipping each second. More time waste is hurting no one else but you, your self. Start creating. Stop criticizing. Start believing. Fear is temporary if it is tamed everyday. It is forever, if it is not tamed everyday. 
h3: <span class=lavendar>3. </span>Sleep More
img.width-half: sleep.jpg
p.text-center.caption: Photo by <a href=https://unsplash.com/@thedakotacorbin target=_blank> Dakota Corbin </a>on Unsplash
p: Please take care of your sleep. 8 hours is best. 6 hours is minimum. Everyday. During day. During night. When ever you get time. Sleep. It is the most important ingredient to creating the most important spaces in your life. You need spiritual lift to be productive. Fresh brain is your best friend. It will fix all of your life’s worries. You will feel better. 
p.note: I was feeling quite dull 2 hours back. Now it is all fine and fresh. 1 hour of sleep has brought me back to life. 
p: Sleep is the final step to achieving the highest place you dream of everyday. Not hardwork. Not timeless efforts. Not boring work routine. Only sleep. 
h3: Conclusion
img: balloons.jpg
p.caption.text-center: Photo by <a href=https://unsplash.com/@artbyhybrid target=_blank> Hybrid </a>on Unsplash
p: Your life is pink, yellow, orange and red. It is not old neither dull. It is lively and colorful. Your 60 years of age doesn’t mean it went out of your hands. Your 20 years doesn’t mean, it is so long that you will bring in colors next year. These colors are for all of us. We all want to live a smiling, chilling and fun day. It is just that we think too much. 
p: Stop thinking. Stop being logical. Stop being all wisdom and correct. Shut off your brain. You are running late. This is time to jump. To tame your fear. To fall. And to rise to success. 
p.note: Guys, if you like this project. Please follow this project's page on twitter. <a href="https://twitter.com/zakatlists">Click here to go to the twitter page</a>. Subscribe below to receive updates.
'@
$ws.Range("E72").Value = "Qasim Ali"

$ws.Range("D72").Select()

"done"
